$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new daily price record was inserted ahead of the existing history,
# pushing every subsequent record (old rows 80-172) down by one row
# (new rows 81-173). Insert a fresh row at row 80 to achieve that shift,
# then populate it with the new record's data.
$ws.Rows.Item(80).Insert()

$ws.Cells.Item(80, 1).Value = 5
$ws.Cells.Item(80, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(80, 3).Value = "Maule"
$ws.Cells.Item(80, 4).Value = 44587
$ws.Cells.Item(80, 5).Value = 7
$ws.Cells.Item(80, 6).Value = 100112021
$ws.Cells.Item(80, 7).Value = "Ají"
$ws.Cells.Item(80, 8).Value = "Americana (o)"
$ws.Cells.Item(80, 9).Value = "Primera"
$ws.Cells.Item(80, 10).Value = 150
$ws.Cells.Item(80, 11).Value = 18000
$ws.Cells.Item(80, 12).Value = 18000
$ws.Cells.Item(80, 13).Value = 18000
$ws.Cells.Item(80, 14).Value = "`$/saco 25 kilos"
$ws.Cells.Item(80, 15).Value = "Región del Maule"
$ws.Cells.Item(80, 16).Value = 720
$ws.Cells.Item(80, 17).Value = 25
$ws.Cells.Item(80, 18).Value = "Hortaliza"
